# Quarterly indexing esoteric bug-fix operation
#
# Column A on the active sheet holds the "as-of" date for each quarterly
# forecast vintage, stored as the 1st of the (quarterly) reference month.
# The correct quarterly index date is the 15th of the month AFTER that
# reference month (i.e. the date was off by roughly one and a half months).
# Re-stamp every populated date in column A (rows 2..last) accordingly,
# leaving the header row (row 1, the forecast-horizon dates across B:BA)
# and every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2

    if ($serial -eq $null -or $serial -eq "") {
        continue
    }

    $d = [DateTime]::FromOADate($serial)
    $d = $d.AddMonths(1)
    $d = $d.AddDays(15 - $d.Day)

    $cell.Value2 = $d.ToOADate()
}
